$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518, shifting existing rows 518-581 down to 519-582
$ws.Rows("518:518").Insert()

# Populate the newly inserted row 518 with the new weekly price record
$ws.Range("A518").Value = 3
$ws.Range("B518").Value = "Femacal de La Calera"
$ws.Range("C518").Value = "Coquimbo"
$ws.Range("D518").Value = 45124
$ws.Range("E518").Value = 5
$ws.Range("F518").Value = 100112009
$ws.Range("G518").Value = "Acelga"
$ws.Range("H518").Value = "Sin especificar"
$ws.Range("I518").Value = "Primera"
$ws.Range("J518").Value = 260
$ws.Range("K518").Value = 3300
$ws.Range("L518").Value = 3500
$ws.Range("M518").Value = 3415
$ws.Range("N518").Value = "$/docena de atados (6 kilos)"
$ws.Range("O518").Value = "Provincia de Quillota"
$ws.Range("P518").Value = 569
$ws.Range("Q518").Value = 6
$ws.Range("R518").Value = "Hortaliza"
